$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Poroto granado" at the
# "Macroferia Regional de Talca" market. It belongs chronologically
# before the existing row 152 entry, so duplicate that row (keeping all
# its formatting/values) and insert the copy above it, shifting rows
# 152:196 down to 153:197.
$ws.Rows.Item(152).Copy()
$ws.Rows.Item(152).Insert()

# Now overwrite the newly inserted row 152 with the new record's data:
# date, min/max/avg price and the per-kilo price differ from the row
# that was duplicated; the remaining columns stay the same.
$ws.Cells.Item(152, 4).Value = 44985   # D152 Fecha
$ws.Cells.Item(152, 11).Value = 30000  # K152 Precio minimo
$ws.Cells.Item(152, 12).Value = 30000  # L152 Precio maximo
$ws.Cells.Item(152, 13).Value = 30000  # M152 Precio promedio ponderado
$ws.Cells.Item(152, 16).Value = 1200   # P152 Precio $/Kg
